# Add "panacea" test data to the Species sheet.
#
# The existing "pugnax" row (row 3) is pushed down to become the new last
# row (row 5), and its old spot (row 3) is filled in with new data for a
# "panacea" test species. The "pugilator" row (row 4) is left untouched.

$wb = $excel.ActiveWorkbook
$wsTraits  = $wb.Worksheets.Item("Traits")
$wsSpecies = $wb.Worksheets.Item("Species")

# --- Species sheet: move the old "pugnax" row down to row 5 ---------------
$wsSpecies.Range("A5").Value = "pugnax"
$wsSpecies.Range("B5").Value = 1.2
$wsSpecies.Range("C5").Value = 2.1
$wsSpecies.Range("D5").Value = 3.1
$wsSpecies.Range("E5").Value = 4.0999999999999996
$wsSpecies.Range("F5").Value = 5.2

# --- Species sheet: fill row 3 with the new "panacea" test data -----------
$wsSpecies.Range("A3").Value = "panacea"
$wsSpecies.Range("B3").Value = 1.2
$wsSpecies.Range("C3").Value = 2.2
$wsSpecies.Range("D3").Value = 3.3
$wsSpecies.Range("E3").Value = 4.3
$wsSpecies.Range("F3").Value = 5.2

# --- Update view/selection state -------------------------------------------
$wsTraits.Activate()
$wsTraits.Range("D11").Select()

$wsSpecies.Activate()
$wsSpecies.Range("A1:F5").Select()
